$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix C32 alignment to match the rest of the US block (center-aligned style)
$ws.Cells.Item(32, 3).HorizontalAlignment = -4108

$xlCenter = -4108
$newIPs = @(
    "47.88.15.127",
    "129.159.84.71",
    "152.67.231.219",
    "193.122.197.111",
    "150.136.61.148",
    "192.9.159.65",
    "47.90.141.204",
    "129.80.116.250",
    "129.146.254.39",
    "192.18.143.199",
    "47.253.56.77",
    "155.248.196.123",
    "192.9.250.241",
    "132.145.134.230",
    "192.9.138.241",
    "47.253.105.131",
    "47.89.244.253",
    "150.230.47.17",
    "164.152.17.14",
    "129.146.243.241",
    "152.70.155.147",
    "129.146.248.140",
    "132.145.152.194",
)

$startRow = 33
for ($i = 0; $i -lt $newIPs.Count; $i++) {
    $r = $startRow + $i
    $ip = $newIPs[$i]

    $ws.Cells.Item($r, 1).Value = $ip

    $ws.Cells.Item($r, 2).Value = ":"
    $ws.Cells.Item($r, 2).HorizontalAlignment = $xlCenter

    $ws.Cells.Item($r, 3).Value = 443
    $ws.Cells.Item($r, 3).HorizontalAlignment = $xlCenter

    $ws.Cells.Item($r, 4).Value = "#"
    $ws.Cells.Item($r, 4).HorizontalAlignment = $xlCenter

    $ws.Cells.Item($r, 5).Value = "US"
    $ws.Cells.Item($r, 5).HorizontalAlignment = $xlCenter

    $ws.Cells.Item($r, 6).Formula = "=ROW()-1"
    $ws.Cells.Item($r, 6).HorizontalAlignment = $xlCenter

    $ws.Cells.Item($r, 7).Formula = "=A" + $r + "&B" + $r + "&C" + $r + "&D" + $r + "&E" + $r + "&""_""&F" + $r
}

# Update view selection to match the extended data range
$ws.Range("G2:G55").Select()
